# Update column F (想去人数 / interest count) values to match the newer scrape snapshot.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 7744
$ws.Range("F9").Value = 8568
$ws.Range("F16").Value = 1170
$ws.Range("F19").Value = 410
$ws.Range("F21").Value = 37
$ws.Range("F22").Value = 564
$ws.Range("F23").Value = 3707
$ws.Range("F24").Value = 73
$ws.Range("F25").Value = 46
$ws.Range("F26").Value = 36
$ws.Range("F28").Value = 3213
$ws.Range("F32").Value = 365
$ws.Range("F33").Value = 142
$ws.Range("F34").Value = 344
$ws.Range("F35").Value = 1019
$ws.Range("F36").Value = 680
$ws.Range("F39").Value = 2685
$ws.Range("F40").Value = 53
$ws.Range("F43").Value = 3224
$ws.Range("F47").Value = 33

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 112
$ws.Range("F3").Value = 138

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 7744
$ws.Range("F9").Value = 8568
$ws.Range("F15").Value = 1170
$ws.Range("F19").Value = 112
$ws.Range("F20").Value = 37
$ws.Range("F21").Value = 138
$ws.Range("F22").Value = 564
$ws.Range("F24").Value = 3707
$ws.Range("F25").Value = 73
$ws.Range("F26").Value = 46
$ws.Range("F27").Value = 36
$ws.Range("F29").Value = 3213
$ws.Range("F31").Value = 365
$ws.Range("F32").Value = 142
$ws.Range("F33").Value = 344
$ws.Range("F35").Value = 1019
$ws.Range("F36").Value = 680
$ws.Range("F40").Value = 2685
$ws.Range("F41").Value = 53
$ws.Range("F44").Value = 3224
$ws.Range("F47").Value = 33

